# trafo_id -> gridnode_id refactor
# The sheet's "gridnode_id" row (A2) held the example value "T1" (a
# transformer/"trafo" id left over from a previous column name). Now that
# the column is properly named gridnode_id, update the sample value to
# reflect the new id scheme: "T0".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("gridnodes")

$ws.Range("A2").Value = "T0"

# Move the active cell/selection down one row, matching the author's
# cursor position when they saved (C7 -> C8).
$ws.Range("C8").Select()
